$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-45 down to 40-46
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with data
$ws.Cells.Item(39, 1).Value = 5
$ws.Cells.Item(39, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(39, 3).Value = "Maule"
$ws.Cells.Item(39, 4).Value = 44449
$ws.Cells.Item(39, 5).Value = 7
$ws.Cells.Item(39, 6).Value = 100112013
$ws.Cells.Item(39, 7).Value = "Alcachofa"
$ws.Cells.Item(39, 8).Value = "Madrigal"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 300
$ws.Cells.Item(39, 11).Value = 12000
$ws.Cells.Item(39, 12).Value = 12000
$ws.Cells.Item(39, 13).Value = 12000
$ws.Cells.Item(39, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(39, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(39, 16).Value = 400
$ws.Cells.Item(39, 17).Value = 30
$ws.Cells.Item(39, 18).Value = "Hortaliza"
